# Updates cryptos list prices/volumes (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.528.65'
$ws.Range("E2").Value = '  +1.32%  '
$ws.Range("D3").Value = '1.881.54'
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4757'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.75%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2922'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06519'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.94'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '97.97'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07722'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7381'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +9.18%  '
$ws.Range("D14").Value = '1.885.61'
$ws.Range("E14").Value = '  +1.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.142'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '273.89'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.71%  '
$ws.Range("D17").Value = '30.528.89'
$ws.Range("E17").Value = '  +1.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.46'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007559'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("D21").Value = '2.132.68'
$ws.Range("E21").Value = '  +1.57%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.246'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.194'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.283'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.942'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1006'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.91%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("E31").Value = '  +5.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.323'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.104'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04815'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.129'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7008'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.13%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01870'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.752'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.312'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.995'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4214'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8439'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.17%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.72'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.340'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.092'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '915.60'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3889'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.34%  '
